$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gewichtung")

# Remove the "barrierefreiheit" criterion row that belonged to the
# "multifunktionale_nutzungsqualitaet" group (was row 10). Deleting the
# whole row shifts the following rows (former 11-13) up by one and drops
# the now-unused shared strings automatically.
$ws.Rows.Item(10).Delete()

# With only 3 remaining criteria in "multifunktionale_nutzungsqualitaet"
# (versickerung, oberflaechentemperatur, befahrbarkeit) their
# within_group_weight becomes 1/3 each instead of 0.25, mirroring the
# existing 1/3 pattern already used for "kreislauffaehigkeit".
$ws.Range("F7").Formula = "=1/3"
$ws.Range("F8:F9").Formula = "=1/3"

# The "kreislauffaehigkeit" group (now rows 10-12 after the deletion)
# keeps its 1/3 weights, rebuilt as a shared formula across its own range.
$ws.Range("F10").Formula = "=1/3"
$ws.Range("F11:F12").Formula = "=1/3"

# Restore the active selection shown in the saved file.
$ws.Activate()
$ws.Range("F18").Select()
